$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" column C for rows 2..28 from 45497 to 45498
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45498
}

# Ensure row 28 has explicit height of 15 (customHeight)
$ws.Rows.Item(28).RowHeight = 15

# Add new row 29 with data for case A 30806-2024
$ws.Cells.Item(29, 1).Value = "A 30806-2024"
$ws.Cells.Item(29, 2).Value = 45497
$ws.Cells.Item(29, 3).Value = 45498
$ws.Cells.Item(29, 4).Value = "OKÄNT"
$ws.Cells.Item(29, 5).Value = "OKÄNT"
$ws.Cells.Item(29, 7).Value = 2.9
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(29, 17).Value = 0

# Apply the same styles as row 28: date format for B/C, wrap-text style for R
$ws.Range("B29:C29").NumberFormat = $ws.Range("B28:C28").NumberFormat
$ws.Cells.Item(29, 18).WrapText = $true

$ws.Rows.Item(29).RowHeight = 15
